# Update countries & provincias Spain
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# --- Swap "Montserrat" / "Islas Malvinas" rows (re-sorted alphabetically) ---
# Row 213 used to be Montserrat (D=12,H=1); Row 214 used to be Islas Malvinas (D=13,H=0).
# After the update they trade places: row 213 -> Islas Malvinas, row 214 -> Montserrat,
# and their D/H figures swap along with the name.
$ws.Range("A213").Value = "Islas Malvinas"
$ws.Range("D213").Value = 13
$ws.Range("H213").Value = 0

$ws.Range("A214").Value = "Montserrat"
$ws.Range("D214").Value = 12
$ws.Range("H214").Value = 1

# --- Refresh the "last updated" banner text ---
$ws.Range("A1").Value = "Datos actualizados a 12 de Agosto de 2020 a las 15:54"

# --- Refresh per-country statistics (columns B..H) ---
# Estados Unidos (row 4)
$ws.Range("B4").Value = 5309622
$ws.Range("C4").Value = 3665
$ws.Range("D4").Value = 2757410
$ws.Range("E4").Value = 2384423
$ws.Range("G4").Value = 40
$ws.Range("H4").Value = 167789

# India (row 6)
$ws.Range("B6").Value = 2360358
$ws.Range("C6").Value = 31953
$ws.Range("D6").Value = 1662457
$ws.Range("E6").Value = 651365
$ws.Range("G6").Value = 348
$ws.Range("H6").Value = 46536

# Arabia Saudita (row 16)
$ws.Range("B16").Value = 293037
$ws.Range("C16").Value = 1569
$ws.Range("D16").Value = 257269
$ws.Range("E16").Value = 32499
$ws.Range("G16").Value = 36
$ws.Range("H16").Value = 3269

# Alemania (row 22)
$ws.Range("B22").Value = 219581
$ws.Range("C22").Value = 51
$ws.Range("E22").Value = 10413

# Irak (row 24)
$ws.Range("B24").Value = 160436
$ws.Range("C24").Value = 3441
$ws.Range("D24").Value = 114541
$ws.Range("E24").Value = 40307
$ws.Range("G24").Value = 57
$ws.Range("H24").Value = 5588

# Catar (row 28)
$ws.Range("B28").Value = 113938
$ws.Range("C28").Value = 292
$ws.Range("D28").Value = 110627
$ws.Range("E28").Value = 3121
$ws.Range("G28").Value = 2
$ws.Range("H28").Value = 190

# Suecia (row 36)
$ws.Range("B36").Value = 83455
$ws.Range("G36").Value = 5
$ws.Range("H36").Value = 5774

# Kuwait (row 41)
$ws.Range("B41").Value = 73785
$ws.Range("C41").Value = 717
$ws.Range("D41").Value = 65451
$ws.Range("E41").Value = 7845
$ws.Range("G41").Value = 3
$ws.Range("H41").Value = 489

# Paises Bajos (row 45)
$ws.Range("B45").Value = 60627
$ws.Range("C45").Value = 654
$ws.Range("G45").Value = 2
$ws.Range("H45").Value = 6161

# Portugal (row 49)
$ws.Range("B49").Value = 53223
$ws.Range("C49").Value = 278
$ws.Range("D49").Value = 38940
$ws.Range("E49").Value = 12519
$ws.Range("G49").Value = 3
$ws.Range("H49").Value = 1764

# Ghana (row 54)
$ws.Range("B54").Value = 41572
$ws.Range("C54").Value = 168
$ws.Range("D54").Value = 39320
$ws.Range("E54").Value = 2029
$ws.Range("G54").Value = 8
$ws.Range("H54").Value = 223

# Azerbaiyan (row 61)
$ws.Range("B61").Value = 33824
$ws.Range("C61").Value = 93
$ws.Range("D61").Value = 31058
$ws.Range("E61").Value = 2269
$ws.Range("G61").Value = 2
$ws.Range("H61").Value = 497

# Serbia (row 63)
$ws.Range("B63").Value = 28751
$ws.Range("C63").Value = 254
$ws.Range("E63").Value = 9128
$ws.Range("G63").Value = 6
$ws.Range("H63").Value = 658

# Republica de Macedonia (row 83)
$ws.Range("B83").Value = 12217
$ws.Range("C83").Value = 134
$ws.Range("D83").Value = 8487
$ws.Range("E83").Value = 3200
$ws.Range("G83").Value = 1
$ws.Range("H83").Value = 530

# Noruega (row 86)
$ws.Range("B86").Value = 9772
$ws.Range("C86").Value = 22
$ws.Range("E86").Value = 659

# Burkina Faso (row 149)
$ws.Range("B149").Value = 1213
$ws.Range("C149").Value = 2
$ws.Range("D149").Value = 995
$ws.Range("E149").Value = 164
